$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1257
$ws.Range("I39").Value = 66
$ws.Range("K39").Value = 198
$ws.Range("M39").Value = 98
$ws.Range("H40").Value = 11115007
$ws.Range("I40").Value = 3619.75
$ws.Range("J40").Value = 55560556
$ws.Range("K40").Value = 3619.75
$ws.Range("L40").Value = 55560556
$ws.Range("M40").Value = -3444.75
$ws.Range("N40").Value = -55560906
$ws.Range("H62").Value = 5317.6665
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376
$ws.Range("H65").Value = 5317.6665
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880
$ws.Range("H97").Value = 949.3333
$ws.Range("I97").Value = 800
$ws.Range("J97").Value = 1024
$ws.Range("K97").Value = 2400
$ws.Range("L97").Value = 3072
$ws.Range("M97").Value = -1904
$ws.Range("N97").Value = -4064
$ws.Range("H107").Value = 775.4545000000001
$ws.Range("I107").Value = 381.1111
$ws.Range("K107").Value = 381.1111
$ws.Range("M107").Value = 1538.8889
$ws.Range("H112").Value = 70212.92999999999
$ws.Range("J112").Value = 75049.57000000001
$ws.Range("L112").Value = 225148.71
$ws.Range("N112").Value = -227364.71
$ws.Range("H125").Value = 12905097
$ws.Range("I125").Value = 2119431.5
$ws.Range("K125").Value = 19074883.5
$ws.Range("M125").Value = -19072423.5
$ws.Range("H133").Value = 90000
$ws.Range("J133").Value = 90000
$ws.Range("L133").Value = 90000
$ws.Range("N133").Value = -100120
$ws.Range("H134").Value = 120000
$ws.Range("J134").Value = 120000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -130140
$ws.Range("H135").Value = 43478900
$ws.Range("I135").Value = 666.5909
$ws.Range("K135").Value = 5999.3181
$ws.Range("M135").Value = -3464.3181
$ws.Range("H137").Value = 2133.926
$ws.Range("I137").Value = 1571.5883
$ws.Range("K137").Value = 4714.7649
$ws.Range("M137").Value = -2164.7649

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6800.095
$ws.Range("J32").Value = 17166.61
$ws.Range("L32").Value = 17166.61
$ws.Range("N32").Value = -17740.61
$ws.Range("H61").Value = 5047.222
$ws.Range("I61").Value = 4365.625
$ws.Range("K61").Value = 4365.625
$ws.Range("M61").Value = -4153.625
$ws.Range("H74").Value = 5722.875
$ws.Range("I74").Value = 5884.7393
$ws.Range("K74").Value = 5884.7393
$ws.Range("M74").Value = -5010.7393
$ws.Range("H77").Value = 5722.875
$ws.Range("I77").Value = 5884.7393
$ws.Range("K77").Value = 29423.6965
$ws.Range("M77").Value = -25055.6965
$ws.Range("H122").Value = 1740.7
$ws.Range("I122").Value = 1972.4286
$ws.Range("K122").Value = 5917.2858
$ws.Range("M122").Value = -3467.2858
$ws.Range("H132").Value = 4794.13
$ws.Range("I132").Value = 4538.864
$ws.Range("K132").Value = 13616.592
$ws.Range("M132").Value = -11086.592
$ws.Range("H136").Value = 5047.222
$ws.Range("I136").Value = 4365.625
$ws.Range("K136").Value = 13096.875
$ws.Range("M136").Value = -10546.875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1950
$ws.Range("I20").Value = 1746.6666
$ws.Range("J20").Value = 2178.75
$ws.Range("K20").Value = 1746.6666
$ws.Range("L20").Value = 2178.75
$ws.Range("M20").Value = -1499.6666
$ws.Range("N20").Value = -2672.75
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("H99").Value = 1858.8889
$ws.Range("J99").Value = 2740.25
$ws.Range("L99").Value = 2740.25
$ws.Range("N99").Value = -5736.25
$ws.Range("M63").ClearContents()
$ws.Range("M66").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14681.954
$ws.Range("I31").Value = 6881.2
$ws.Range("J31").Value = 16976.295
$ws.Range("K31").Value = 6881.2
$ws.Range("L31").Value = 16976.295
$ws.Range("M31").Value = -6586.2
$ws.Range("N31").Value = -17566.295
$ws.Range("H34").Value = 14681.954
$ws.Range("I34").Value = 6881.2
$ws.Range("J34").Value = 16976.295
$ws.Range("K34").Value = 6881.2
$ws.Range("L34").Value = 16976.295
$ws.Range("M34").Value = -6679.2
$ws.Range("N34").Value = -17380.295
$ws.Range("H107").Value = 615529.0600000001
$ws.Range("I107").Value = 724868.2
$ws.Range("K107").Value = 724868.2
$ws.Range("M107").Value = -722948.2
$ws.Range("H131").Value = 90096
$ws.Range("I131").Value = 55555
$ws.Range("K131").Value = 55555
$ws.Range("M131").Value = -50515
$ws.Range("H132").Value = 35489.2
$ws.Range("I132").Value = 37881.355
$ws.Range("K132").Value = 113644.065
$ws.Range("M132").Value = -111114.065
$ws.Range("H134").Value = 2502.2
$ws.Range("I134").Value = 1893.7693
$ws.Range("J134").Value = 6457
$ws.Range("K134").Value = 5681.3079
$ws.Range("L134").Value = 19371
$ws.Range("M134").Value = -3146.3079
$ws.Range("N134").Value = -24441

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 11.5
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 11.5
$ws.Range("K2").Value = 0
$ws.Range("N2").Value = -295
$ws.Range("H63").Value = 1600
$ws.Range("I63").Value = 1600
$ws.Range("K63").Value = 4800
$ws.Range("M63").Value = -4051
$ws.Range("H64").Value = 4764.2
$ws.Range("I64").Value = 4764.2
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 14292.6
$ws.Range("L64").Value = 0
$ws.Range("H66").Value = 1600
$ws.Range("I66").Value = 1600
$ws.Range("K66").Value = 14400
$ws.Range("M66").Value = -10656
$ws.Range("H67").Value = 4764.2
$ws.Range("I67").Value = 4764.2
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 14292.6
$ws.Range("L67").Value = 0
$ws.Range("H122").Value = 1661.6
$ws.Range("I122").Value = 922.55554
$ws.Range("J122").Value = 2770.1667
$ws.Range("K122").Value = 8302.99986
$ws.Range("L122").Value = 24931.5003
$ws.Range("M122").Value = -5852.99986
$ws.Range("N122").Value = -29831.5003
$ws.Range("H123").Value = 7724.75
$ws.Range("J123").Value = 8666.333000000001
$ws.Range("L123").Value = 25998.999
$ws.Range("N123").Value = -30898.999
$ws.Range("H131").Value = 1621.8
$ws.Range("J131").Value = 1641.6154
$ws.Range("L131").Value = 4924.8462
$ws.Range("N131").Value = -15004.8462
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 56725
$ws.Range("J109").Value = 56725
$ws.Range("L109").Value = 56725
$ws.Range("N109").Value = -58805
$ws.Range("H122").Value = 82026.266
$ws.Range("I122").Value = 121639.5
$ws.Range("J122").Value = 2799.8
$ws.Range("K122").Value = 364918.5
$ws.Range("L122").Value = 8399.400000000001
$ws.Range("M122").Value = -362468.5
$ws.Range("N122").Value = -13299.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1436.619
$ws.Range("J22").Value = 1040.8334
$ws.Range("L22").Value = 1040.8334
$ws.Range("N22").Value = -1630.8334
$ws.Range("H27").Value = 1436.619
$ws.Range("J27").Value = 1040.8334
$ws.Range("L27").Value = 1040.8334
$ws.Range("N27").Value = -1254.8334
$ws.Range("H46").Value = 1229.3334
$ws.Range("I46").Value = 960
$ws.Range("K46").Value = 960
$ws.Range("M46").Value = -772

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2323.9167
$ws.Range("I126").Value = 2398.818
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 7196.454000000001
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -4726.454000000001
$ws.Range("N126").Value = -9440
$ws.Range("H132").Value = 21367.584
$ws.Range("I132").Value = 21491.908
$ws.Range("K132").Value = 64475.724
$ws.Range("M132").Value = -61945.724

Write-Host "Applied Spriggan_Profits updates"